$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LFHVM")

# Remove the extra Summer/Winter Peak 2-5 rows (rows 8-15), keeping a single
# "Summer Peak" slot (row 6) and a single "Winter Peak" slot (row 7).
$ws.Range("A8:I15").EntireRow.Delete()

# Rename the remaining peak timeslice labels.
$ws.Range("A6").Value = "Summer Peak"
$ws.Range("A7").Value = "Winter Peak"

# Restore the cursor/selection state seen in the saved file.
$ws.Range("F25").Select() | Out-Null
$wb.Worksheets.Item("About").Select() | Out-Null
